# Remove the trailing "a" from the AMU labels in column A (rows 11-20),
# e.g. "AMU32a" -> "AMU32". Other cells / rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 11; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = [string]$cell.Value()
    if ($val.EndsWith("a")) {
        $cell.Value = $val.Substring(0, $val.Length - 1)
    }
}

# Update the active selection as recorded in the saved workbook view.
$ws.Range("A12").Select()
